$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet title (tab name) to reflect the new "through" date
$ws.Name = "Through 2022-03-09"

# Update the column header text for the "through" month label
$ws.Range("B1").Value = "March 2022 (through March 09)"

# Row 3 - Austin
$ws.Range("E3").Value = 2
$ws.Range("Q3").Value = 1

# Row 4 - North Lawndale
$ws.Range("K4").Value = 2

# Row 5 - Garfield Park
$ws.Range("E5").Value = 2
$ws.Range("T5").Value = 3

# Row 12 - Englewood
$ws.Range("H12").Value = 1

# Row 15 - Humboldt Park
$ws.Range("N15").Value = 1

# Row 21 - West Pullman
$ws.Range("T21").Value = 3

# Row 25 - Washington Park
$ws.Range("B25").Value = 1

# Row 26 - Grand Crossing
$ws.Range("E26").Value = 2

# Row 27 - Calumet Heights
$ws.Range("E27").Value = 2

# Row 33 - Belmont Cragin
$ws.Range("B33").Value = 3

# Row 44 - Grand Boulevard
$ws.Range("E44").Value = 1

# Row 49 - Little Village
$ws.Range("B49").Value = 1

# Row 50 - Albany Park
$ws.Range("Q50").Value = 1

# Row 54 - Avalon Park
$ws.Range("N54").Value = 1

# Row 56 - Clearing
$ws.Range("E56").Value = 2
